# Commiting scripts (R22 UAT2 - Regression).
# Append two new Transaction Number rows to the "Outward Clearing" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "OT2325017475"
$ws.Range("A6").Value = "OT2325009833"
